$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill the props table (rows 1-15, cols A-E) ---
$data = @(
  @('Property','Description','Accepts','Default','Example'),
  @('fill','The fill color of the SVG','string','```''inherit''```','```fill="red"```'),
  @('fillOpacity','The opacity of the SVG fill','float','```''inherit''```','```fillOpacity={0.2}```'),
  @('stroke','The stroke color for the SVG','string','```''inherit''```','```stroke="red"```'),
  @('strokeWidth','The width of the SVG stroke','string','```''inherit''```','```strokeWidth="20px"```'),
  @('width','SVG width','string or number','```''inherit''```','```width="20px"```'),
  @('height','SVG height','string ornumber','```''inherit''```','```height="20px"```'),
  @('animation','The name of the animation desired, if only one animation in use','string','```''none''```','```animation=''fade-in-stroke''```'),
  @('duration','The duration of the animation desired, if only one animation in use','string','```''0.5s''```','```duration=''2.5s''```'),
  @('iterationCount','The desired iteration count, if only one animation in use','string or number','```1```','```iterationCount=''infinite''```'),
  @('timingFunction','The desired timing function, if only one animation in use','string','null','```timingFunction=''linear''```'),
  @('animationNames','An array of animation names desired','Array','null','```animationNames={[''fade-in-fill'', ''fade-in-stroke'']}```'),
  @('animationDurations','An array of animation durations (maps 1 to 1 to animation names)','Array','null (if animation name is present without duration, defaults to ''0.5s'')','```animationDurations={[''0.5s'', ''3s'']}```'),
  @('animationIterationCounts','An array of animation iteration counts (maps 1 to 1 to animation names)','Array','null (if animation name is present without duration, defaults to ''1'')','```animationIterationCounts={[''infinite'','''',3]}```'),
  @('animationTimingFunctions','An array of animation timing functions (maps 1 to 1 to animation names)','Array','null (if animation name is present without duration, defaults to none)','```animationTimingFunctions={[''linear'','''','''']}```')
)

for ($r = 0; $r -lt $data.Length; $r++) {
  for ($c = 0; $c -lt $data[$r].Length; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value2 = $data[$r][$c]
  }
}

# --- Column widths (bestFit-style autosized columns A:E) ---
$ws.Columns.Item(1).ColumnWidth = 24.17
$ws.Columns.Item(2).ColumnWidth = 65.67
$ws.Columns.Item(3).ColumnWidth = 13.33
$ws.Columns.Item(4).ColumnWidth = 63.67
$ws.Columns.Item(5).ColumnWidth = 47.5

# --- Selection moves to C17 ---
$null = $ws.Range("C17").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

